$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal TEXT (not auto-converted to a number),
# and reset the cell style back to "Normal" afterwards so no stray style survives.
function Set-TextCell {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Clear out the old data rows (rows 2-16), keeping the header row (row 1) untouched.
$ws.Range("A2:E16").ClearContents()

# New cohort-retention data: (cohort_year, period_index, num_customers, cohort_size, retention_rate)
$data = @(
    @("2020", 0, 109, 109, 1),
    @("2020", 1, 66, 109, 0.6055045871559633),
    @("2020", 2, 27, 109, 0.2477064220183486),
    @("2020", 3, 18, 109, 0.1651376146788991),
    @("2020", 4, 16, 109, 0.1467889908256881),
    @("2020", 5, 14, 109, 0.1284403669724771),
    @("2021", 0, 419, 419, 1),
    @("2021", 1, 142, 419, 0.3389021479713604),
    @("2021", 2, 106, 419, 0.2529832935560859),
    @("2021", 3, 87, 419, 0.20763723150358),
    @("2021", 4, 66, 419, 0.1575178997613365),
    @("2022", 0, 193, 193, 1),
    @("2022", 1, 73, 193, 0.3782383419689119),
    @("2022", 2, 63, 193, 0.3264248704663212),
    @("2022", 3, 57, 193, 0.2953367875647668),
    @("2023", 0, 123, 123, 1),
    @("2023", 1, 85, 123, 0.6910569105691057),
    @("2023", 2, 61, 123, 0.4959349593495935),
    @("2024", 0, 206, 206, 1),
    @("2024", 1, 107, 206, 0.5194174757281553),
    @("2025", 0, 35, 35, 1)
)

$rowNum = 2
foreach ($rec in $data) {
    Set-TextCell $ws.Cells.Item($rowNum, 1) $rec[0]
    $ws.Cells.Item($rowNum, 2).Value = $rec[1]
    $ws.Cells.Item($rowNum, 3).Value = $rec[2]
    $ws.Cells.Item($rowNum, 4).Value = $rec[3]
    $ws.Cells.Item($rowNum, 5).Value = $rec[4]
    $rowNum++
}
